# CA 4.0 files test
# Update the "About" sheet of the Output Currency Conversion Factors workbook:
#  - refresh the currency-year labels from 2021 to 2020
#  - replace the 2012/2021 CPI-adjustment factor with a new 2019->2012 factor
#    (entered as a literal value, no longer the "=1/1.21" formula)
#  - stamp a "date updated" value in C1
#  - leave the selection on B31, matching the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Date stamp in C1 (2022-03-11), formatted as a short date
$ws.Range("C1").Value = 44631
$ws.Range("C1").NumberFormat = "mm-dd-yy"

# Year labels: 2021 -> 2020
$ws.Range("A18").Value = "billion 2020 dollars"
$ws.Range("A21").Value = "million 2020 dollars"
$ws.Range("A24").Value = "2020 dollars"

# Conversion factor: was "=1/1.21" (2012$ per 2021$); now a literal 2019$-per-2012$ factor
$ws.Range("A26").Value = 0.88711067149387013
$ws.Range("B26").Value = "2019 dollars per 2012 dollar"

# Explanatory note below also refers to the conversion year
$ws.Range("B29").Value = "which in this case is ""2012 dollars per 2020 dollar."""

$ws.Range("B31").Select() | Out-Null
